$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell "Source" in column F, row 1
$ws.Range("F1").Value = "Source"

# Hyperlink cell F4 pointing to the wikipedia article (creates the Hyperlink style + font)
$ws.Hyperlinks.Add($ws.Range("F4"), "https://en.wikipedia.org/wiki/List_of_countries_by_system_of_government")

# Bold the new "Source" header to match the rest of row 1
$ws.Range("F1").Font.Bold = $true

# New data row 4 - system_of_gov dataset
$ws.Range("A4").Value = "system_of_gov"
$ws.Range("B4").Value = "wikipedia article of systems of government (democracy, chiefdom, etc.) for all UN countries"

# Update selection to match target state
$ws.Range("A7").Select()
